# Update the "Input" sheet with the new PRO / serial number / cell info.
# Downstream formulas (Input!B7 and everything on Template_printout) will
# recalculate automatically because they reference these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

$ws.Range("B3").Value = 20578713
$ws.Range("B4").Value = "A01593"
$ws.Range("B5").Value = "APXCAS2131001"
